$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are number-like strings (locale-formatted prices using "." as
# thousands separator) that Excel would otherwise auto-convert to numbers/dates.
# Force them to remain plain text by temporarily applying a text number format,
# then restore the default (Normal) style so no residual formatting is left on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.119.90"
$ws.Range("E2").Value = "  +0.14%  "
Set-TextValue $ws.Range("D3") "1.656.42"
$ws.Range("E3").Value = "  -0.24%  "
Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  -0.19%  "
Set-TextValue $ws.Range("D5") "214.90"
$ws.Range("E5").Value = "  +3.27%  "
Set-TextValue $ws.Range("D6") "0.5234"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("E7").Value = "  -0.16%  "
Set-TextValue $ws.Range("D8") "0.2622"
$ws.Range("E8").Value = "  +1.60%  "
Set-TextValue $ws.Range("D9") "0.06390"
$ws.Range("E9").Value = "  +1.61%  "
Set-TextValue $ws.Range("D10") "20.83"
$ws.Range("E10").Value = "  -0.44%  "
Set-TextValue $ws.Range("D11") "0.07753"
$ws.Range("E11").Value = "  +2.89%  "
Set-TextValue $ws.Range("D12") "1.659.54"
$ws.Range("E12").Value = "  -0.18%  "
Set-TextValue $ws.Range("D13") "4.451"
$ws.Range("E13").Value = "  +1.11%  "
Set-TextValue $ws.Range("D14") "1.881.99"
$ws.Range("E14").Value = "  -0.30%  "
Set-TextValue $ws.Range("D15") "0.5512"
$ws.Range("E15").Value = "  +2.23%  "
Set-TextValue $ws.Range("D16") "0.0₅8271"
$ws.Range("E16").Value = "  +4.51%  "
Set-TextValue $ws.Range("D17") "65.13"
$ws.Range("E17").Value = "  -1.47%  "
Set-TextValue $ws.Range("D18") "26.127.38"
$ws.Range("E18").Value = "  +0.16%  "
Set-TextValue $ws.Range("D20") "4.757"
$ws.Range("E20").Value = "  +1.36%  "
Set-TextValue $ws.Range("D21") "190.31"
$ws.Range("E21").Value = "  +1.50%  "
$ws.Range("E22").Value = "  +1.30%  "
Set-TextValue $ws.Range("D23") "6.355"
$ws.Range("E23").Value = "  +2.69%  "
$ws.Range("E24").Value = "  -0.19%  "
Set-TextValue $ws.Range("D25") "143.09"
$ws.Range("E25").Value = "  -3.32%  "
Set-TextValue $ws.Range("D26") "0.1251"
$ws.Range("E26").Value = "  +3.74%  "
Set-TextValue $ws.Range("D27") "7.402"
$ws.Range("E27").Value = "  +0.32%  "
Set-TextValue $ws.Range("D28") "16.00"
$ws.Range("E28").Value = "  +2.42%  "
Set-TextValue $ws.Range("D29") "1.432"
$ws.Range("E29").Value = "  +3.55%  "
Set-TextValue $ws.Range("D30") "0.06041"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("E31").Value = "  +0.17%  "
Set-TextValue $ws.Range("D32") "3.503"
$ws.Range("E32").Value = "  +1.03%  "
Set-TextValue $ws.Range("D33") "3.415"
$ws.Range("E33").Value = "  +0.62%  "
Set-TextValue $ws.Range("D34") "1.658"
$ws.Range("E34").Value = "  +1.53%  "
Set-TextValue $ws.Range("D35") "0.9988"
$ws.Range("E35").Value = "  +1.59%  "
Set-TextValue $ws.Range("D36") "2.400"
Set-TextValue $ws.Range("D37") "2.757"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  -3.91%  "
Set-TextValue $ws.Range("D39") "0.01604"
$ws.Range("E39").Value = "  +0.74%  "
$ws.Range("E40").Value = "  -0.76%  "
Set-TextValue $ws.Range("D41") "0.8539"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("E42").Value = "  -0.21%  "
Set-TextValue $ws.Range("D43") "1.029.20"
$ws.Range("E43").Value = "  -6.74%  "
Set-TextValue $ws.Range("D44") "99.63"
$ws.Range("E44").Value = "  -0.20%  "
Set-TextValue $ws.Range("D45") "1.803.78"
Set-TextValue $ws.Range("D46") "0.0₈106"
$ws.Range("E46").Value = "  -3.16%  "
Set-TextValue $ws.Range("D47") "55.96"
$ws.Range("E47").Value = "  +1.68%  "
Set-TextValue $ws.Range("D48") "1.001"
$ws.Range("E48").Value = "  +0.25%  "
Set-TextValue $ws.Range("D49") "8.105"
$ws.Range("E49").Value = "  +0.98%  "
Set-TextValue $ws.Range("D50") "0.05156"
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D51") "0.4215"
$ws.Range("E51").Value = "  -0.54%  "
